# Adds a new "valorRecarga" column (M) to the data-driven sheet used by the
# "CargarEprepago" task: a header cell (M1) styled like the other green
# header cells but with only a left/right border, and a data cell (M2)
# holding the recharge amount 10000 (stored as text-formatted number, like
# the rest of the data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- M1: header cell "valorRecarga" ---------------------------------------
# Reuse the exact header look (bold "Mic Shell Dlg" font + green fill) from
# an existing header cell so no new font/fill is introduced, then narrow the
# border to left/right only.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("M1").Value = "valorRecarga"
$ws.Range("M1").Borders.Item(8).LineStyle = -4142   # xlEdgeTop    -> xlLineStyleNone
$ws.Range("M1").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> xlLineStyleNone

# --- M2: data cell with the recharge amount --------------------------------
$ws.Range("M2").Value = 10000
$ws.Range("M2").NumberFormat = "@"

# Match the saved selection state (active cell on the new column).
$ws.Range("M2").Select() | Out-Null
